$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textCells = @("D5", "D6", "D9", "D10", "D13", "D15", "D16", "D18", "D21", "D22", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D35", "D37", "D38", "D40", "D41", "D43", "D47", "D48", "D50", "D51")
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range("D2").Value = "26.331.80"
$ws.Range("E2").Value = "  +1.03%  "
$ws.Range("D3").Value = "1.680.93"
$ws.Range("E3").Value = "  +0.88%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "218.26"
$ws.Range("E5").Value = "  +0.80%  "
$ws.Range("D6").Value = "0.5282"
$ws.Range("E6").Value = "  +3.51%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("E8").Value = "  +2.67%  "
$ws.Range("D9").Value = "0.06470"
$ws.Range("D10").Value = "21.98"
$ws.Range("E10").Value = "  +1.47%  "
$ws.Range("E11").Value = "  +1.39%  "
$ws.Range("D12").Value = "1.688.80"
$ws.Range("E12").Value = "  +1.32%  "
$ws.Range("D13").Value = "4.518"
$ws.Range("E13").Value = "  +0.31%  "
$ws.Range("E14").Value = "  -0.27%  "
$ws.Range("D15").Value = "0.000008517"
$ws.Range("E15").Value = "  -0.27%  "
$ws.Range("D16").Value = "64.85"
$ws.Range("E16").Value = "  +0.83%  "
$ws.Range("D17").Value = "26.350.40"
$ws.Range("E17").Value = "  +0.82%  "
$ws.Range("D18").Value = "4.924"
$ws.Range("E18").Value = "  +0.27%  "
$ws.Range("E19").Value = "  +0.16%  "
$ws.Range("E20").Value = "  +1.11%  "
$ws.Range("D21").Value = "190.13"
$ws.Range("E21").Value = "  +0.91%  "
$ws.Range("D22").Value = "6.210"
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("D24").Value = "144.88"
$ws.Range("E24").Value = "  -0.52%  "
$ws.Range("D25").Value = "7.784"
$ws.Range("E25").Value = "  +2.14%  "
$ws.Range("D26").Value = "0.1257"
$ws.Range("E26").Value = "  +5.64%  "
$ws.Range("D27").Value = "15.81"
$ws.Range("E27").Value = "  +1.48%  "
$ws.Range("D28").Value = "0.06515"
$ws.Range("E28").Value = "  +1.07%  "
$ws.Range("D29").Value = "1.363"
$ws.Range("E29").Value = "  +4.51%  "
$ws.Range("D30").Value = "1.326"
$ws.Range("E30").Value = "  +0.65%  "
$ws.Range("D31").Value = "3.593"
$ws.Range("E31").Value = "  +2.00%  "
$ws.Range("D32").Value = "3.591"
$ws.Range("E32").Value = "  +2.34%  "
$ws.Range("E33").Value = "  +1.74%  "
$ws.Range("E34").Value = "  +1.22%  "
$ws.Range("D35").Value = "0.6221"
$ws.Range("E35").Value = "  +2.77%  "
$ws.Range("E36").Value = "  +1.81%  "
$ws.Range("D37").Value = "2.738"
$ws.Range("E37").Value = "  +1.91%  "
$ws.Range("D38").Value = "6.300"
$ws.Range("E38").Value = "  +1.78%  "
$ws.Range("D39").Value = "1.116.11"
$ws.Range("E39").Value = "  +3.89%  "
$ws.Range("D40").Value = "0.01622"
$ws.Range("E40").Value = "  +0.63%  "
$ws.Range("D41").Value = "0.8747"
$ws.Range("E41").Value = "  +1.80%  "
$ws.Range("E42").Value = "  +0.68%  "
$ws.Range("D43").Value = "100.55"
$ws.Range("E43").Value = "  -0.16%  "
$ws.Range("D44").Value = "1.829.85"
$ws.Range("E44").Value = "  +0.83%  "
$ws.Range("E45").Value = "  -2.17%  "
$ws.Range("E46").Value = "  +1.46%  "
$ws.Range("D47").Value = "8.173"
$ws.Range("E47").Value = "  +1.64%  "
$ws.Range("D48").Value = "1.003"
$ws.Range("E48").Value = "  -0.17%  "
$ws.Range("E49").Value = "  +1.21%  "
$ws.Range("D50").Value = "0.4296"
$ws.Range("E50").Value = "  -0.03%  "
$ws.Range("D51").Value = "6.090"
$ws.Range("E51").Value = "  +2.38%  "
